$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = "2022/2023, 2025/2026"
    22 = "2024/2025, 2025/2026"
    23 = "2022/2023, 2025/2026, 2023/2024"
    24 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    27 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    28 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    31 = "2022/2023, 2025/2026"
    50 = "2024/2025, 2025/2026"
    51 = "2022/2023, 2025/2026, 2023/2024"
    52 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    55 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    56 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
